{"js": "// Update the two-digit multiplication problems in the table.\n// Each \"AxB=\" string in the document is unique, so a simple\n// search + insertText(\"Replace\") pass for each (old, new) pair\n// reproduces the diff.\n\nconst pairs = [\n  [\"90\u00d743=\", \"70\u00d723=\"],\n  [\"93\u00d794=\", \"89\u00d755=\"],\n  [\"71\u00d745=\", \"37\u00d799=\"],\n  [\"43\u00d789=\", \"73\u00d721=\"],\n  [\"52\u00d713=\", \"56\u00d733=\"],\n  [\"16\u00d751=\", \"21\u00d758=\"],\n  [\"43\u00d717=\", \"20\u00d789=\"],\n  [\"53\u00d730=\", \"83\u00d721=\"],\n  [\"28\u00d752=\", \"61\u00d737=\"],\n  [\"24\u00d753=\", \"28\u00d751=\"],\n  [\"13\u00d737=\", \"82\u00d736=\"],\n  [\"95\u00d782=\", \"31\u00d780=\"],\n  [\"21\u00d798=\", \"39\u00d785=\"],\n  [\"75\u00d793=\", \"74\u00d790=\"],\n  [\"60\u00d712=\", \"68\u00d784=\"],\n  [\"20\u00d762=\", \"82\u00d778=\"],\n  [\"99\u00d784=\", \"89\u00d734=\"],\n  [\"18\u00d718=\", \"75\u00d781=\"],\n  [\"74\u00d791=\", \"20\u00d738=\"],\n  [\"17\u00d776=\", \"64\u00d743=\"],\n  [\"62\u00d758=\", \"64\u00d718=\"],\n  [\"55\u00d724=\", \"12\u00d733=\"],\n  [\"11\u00d743=\", \"84\u00d769=\"],\n  [\"47\u00d759=\", \"12\u00d770=\"],\n  [\"77\u00d713=\", \"86\u00d795=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the two-digit multiplication problems in the table.\n# Each \"AxB=\" string in the document is unique, so a simple\n# Find/Replace pass for each (old, new) pair reproduces the diff.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"90\u00d743=\", \"70\u00d723=\"),\n    @(\"93\u00d794=\", \"89\u00d755=\"),\n    @(\"71\u00d745=\", \"37\u00d799=\"),\n    @(\"43\u00d789=\", \"73\u00d721=\"),\n    @(\"52\u00d713=\", \"56\u00d733=\"),\n    @(\"16\u00d751=\", \"21\u00d758=\"),\n    @(\"43\u00d717=\", \"20\u00d789=\"),\n    @(\"53\u00d730=\", \"83\u00d721=\"),\n    @(\"28\u00d752=\", \"61\u00d737=\"),\n    @(\"24\u00d753=\", \"28\u00d751=\"),\n    @(\"13\u00d737=\", \"82\u00d736=\"),\n    @(\"95\u00d782=\", \"31\u00d780=\"),\n    @(\"21\u00d798=\", \"39\u00d785=\"),\n    @(\"75\u00d793=\", \"74\u00d790=\"),\n    @(\"60\u00d712=\", \"68\u00d784=\"),\n    @(\"20\u00d762=\", \"82\u00d778=\"),\n    @(\"99\u00d784=\", \"89\u00d734=\"),\n    @(\"18\u00d718=\", \"75\u00d781=\"),\n    @(\"74\u00d791=\", \"20\u00d738=\"),\n    @(\"17\u00d776=\", \"64\u00d743=\"),\n    @(\"62\u00d758=\", \"64\u00d718=\"),\n    @(\"55\u00d724=\", \"12\u00d733=\"),\n    @(\"11\u00d743=\", \"84\u00d769=\"),\n    @(\"47\u00d759=\", \"12\u00d770=\"),\n    @(\"77\u00d713=\", \"86\u00d795=\")\n)\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair[0]\n    $find.Replacement.Text = $pair[1]\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
